$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three rows whose content changes (old ElasticSearch/ReactNative/SmartTube rows
# become javascript-algorithms / javascript_algorithms_2 / QuixBugs(Java_Corrected))
$ws.Range("A19").Value = "javascript-algorithms"
$ws.Range("A20").Value = "javascript_algorithms_2"
$ws.Range("A21").Value = "QuixBugs(Java_Corrected)"

# Remove the now-duplicate trailing rows 22-24 (old javascript-algorithms,
# javascript-algorithms-2, QuixBugs(Java_Corrected) entries) by deleting the rows
# entirely so the used range shrinks back down.
$ws.Rows("22:24").Delete()

# Move the selection back to A14, matching the saved view state in the target file.
$ws.Range("A14").Select()
